# checkpoint da automação das figuras do grupo 13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row fix: typo "Trimestre" -> "Trimesetre" ---
$ws.Range("C1").Value = "Trimesetre"

# --- Variable label: drop the "Taxa de " prefix for rows 2-8 ---
$novoLabel = "Pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência"
$ws.Range("B2").Value = $novoLabel
$ws.Range("B3").Value = $novoLabel
$ws.Range("B4").Value = $novoLabel
$ws.Range("B5").Value = $novoLabel
$ws.Range("B6").Value = $novoLabel
$ws.Range("B7").Value = $novoLabel
$ws.Range("B8").Value = $novoLabel

# --- Value column: round to 2 decimals ---
$ws.Range("D2").Value = 56.37
$ws.Range("D3").Value = 55.26
$ws.Range("D4").Value = 54.5
$ws.Range("D5").Value = 54.48
$ws.Range("D6").Value = 54.42
$ws.Range("D7").Value = 54.28
$ws.Range("D8").Value = 46.79

# --- Rows 9/10: swap Nordeste/Brasil and refresh the value + label/date ---
# (the date text is copied from a cell that already holds it as text, so the
#  "01/07/2024" stays a string instead of being reinterpreted as a date)
$ws.Range("A9").Value = "Brasil"
$ws.Range("B9").Value = $novoLabel
$ws.Range("C2").Copy($ws.Range("C9"))
$ws.Range("D9").Value = 50.67
$ws.Range("E9").Value = " "

$ws.Range("A10").Value = "Nordeste"
$ws.Range("B10").Value = $novoLabel
$ws.Range("C2").Copy($ws.Range("C10"))
$ws.Range("D10").Value = 43.92
$ws.Range("E10").Value = " "

# --- Page margins reset to Excel's defaults (inches -> points, 1in = 72pt) ---
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72
